# Apply update to t11.4 sheet: column B changes from text quarter labels
# (e.g. "2010 4 trim") to plain numeric years, and column C values are
# refreshed with updated figures. Header row (B1/C1) keeps its "Data"/"Valor"
# labels untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years = @(2010,2010,2010,2010,2011,2011,2011,2011,2012,2012,2012,2012,2013,2013,2013,2013,2014,2014,2014,2014,2015,2015,2015,2015,2016,2016,2016,2016,2017,2017,2017,2017,2018,2018,2018,2018,2019,2019,2019,2019,2020,2020,2020,2020,2021,2021,2021,2021,2022,2022,2022,2022,2023,2023,2023,2023)
$values = @(3098.5809705903548,3270.779696923239,23.24947726611369,-195.46939728839249,3822.763517699701,3911.1772244567001,12.735712879693754,-101.1494196366927,4851.6141492202578,4773.9627523969748,74.755921009532415,2.8954758137494951,4597.186655099481,4598.5535973779897,4.1008268355266457,-5.4677691140355273,4875.7880913134968,4671.1648862029751,172.95791010034841,31.665295010172791,6050.4584366048948,5821.5356978258142,229.25437978086754,-0.33164100178715711,7444.0303730971409,7341.3045012668317,155.88292935067116,-53.171240299851355,6601.6396040240334,6272.3484076347613,328.85034130248181,0.4270783653278985,7028.0738465203176,6869.8373114392198,236.71739462534867,-78.480859544250578,4432.5531621446362,4565.325931901697,-43.082083887035388,-89.67795475776515,6294.3677912056082,6335.5268756372961,127.50911388899378,-168.66819832068174,6036.5046240057254,5890.4154415820813,140.05746996751657,6.0427798184325594,6186.2823513574858,5792.5509648923617,266.5538678949473,127.16705645117689,6470.48,5831.48,551.28,87.72)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $years[$i]
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
